$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("K8").Value = 1500
$ws.Range("M8").Value = -1361

$ws.Range("H43").Value = 5107.143
$ws.Range("J43").Value = 5291.6665
$ws.Range("L43").Value = 5291.6665
$ws.Range("N43").Value = -5429.6665

$ws.Range("H55").Value = 762.8421
$ws.Range("I55").Value = 388.8889
$ws.Range("J55").Value = 1099.4
$ws.Range("K55").Value = 388.8889
$ws.Range("L55").Value = 1099.4
$ws.Range("M55").Value = -174.8889
$ws.Range("N55").Value = -1527.4

$ws.Range("H64").Value = 6498.7646
$ws.Range("J64").Value = 8222.223
$ws.Range("L64").Value = 8222.223
$ws.Range("N64").Value = -8718.223

$ws.Range("H67").Value = 6498.7646
$ws.Range("J67").Value = 8222.223
$ws.Range("L67").Value = 8222.223
$ws.Range("N67").Value = -9938.223

$ws.Range("H108").Value = 63207
$ws.Range("J108").Value = 63207
$ws.Range("L108").Value = 63207
$ws.Range("N108").Value = -70887

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H111").Value = 1873.8334
$ws.Range("J111").Value = 2097.6
$ws.Range("L111").Value = 6292.799999999999
$ws.Range("N111").Value = -12426.8

$ws.Range("H116").Value = 12491.333
$ws.Range("I116").Value = 14006.706
$ws.Range("J116").Value = 6051
$ws.Range("K116").Value = 14006.706
$ws.Range("L116").Value = 6051
$ws.Range("M116").Value = -10564.706
$ws.Range("N116").Value = -12935

$ws.Range("H132").Value = 1402.5077
$ws.Range("I132").Value = 1401.875
$ws.Range("K132").Value = 4205.625
$ws.Range("M132").Value = -1675.625

$ws.Range("H137").Value = 1917680.1
$ws.Range("I137").Value = 1824.5333
$ws.Range("K137").Value = 5473.5999
$ws.Range("M137").Value = -2923.5999

$ws.Range("H138").Value = 2420.02
$ws.Range("I138").Value = 1223.3684
$ws.Range("J138").Value = 2700.716
$ws.Range("K138").Value = 3670.1052
$ws.Range("L138").Value = 8102.147999999999
$ws.Range("M138").Value = 1469.8948
$ws.Range("N138").Value = -18382.148

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H5").Value = 999
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H20").Value = 37764.43
$ws.Range("I20").Value = 56966.89
$ws.Range("J20").Value = 3200
$ws.Range("K20").Value = 56966.89
$ws.Range("L20").Value = 3200
$ws.Range("M20").Value = -56719.89
$ws.Range("N20").Value = -3694

$ws.Range("H107").Value = 1331.6666
$ws.Range("I107").Value = 998.5
$ws.Range("K107").Value = 998.5
$ws.Range("M107").Value = 921.5

$ws.Range("H134").Value = 2860556.5
$ws.Range("I134").Value = 3573750
$ws.Range("J134").Value = 7782.6
$ws.Range("K134").Value = 10721250
$ws.Range("L134").Value = 23347.8
$ws.Range("M134").Value = -10718715
$ws.Range("N134").Value = -28417.8

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 4876.8
$ws.Range("I31").Value = 2431.5
$ws.Range("K31").Value = 2431.5
$ws.Range("M31").Value = -2136.5

$ws.Range("H34").Value = 4876.8
$ws.Range("I34").Value = 2431.5
$ws.Range("K34").Value = 2431.5
$ws.Range("M34").Value = -2229.5

$ws.Range("H100").Value = 60006.332
$ws.Range("J100").Value = 60006.332
$ws.Range("L100").Value = 60006.332
$ws.Range("N100").Value = -62170.332

$ws.Range("H120").Value = 69999.5
$ws.Range("J120").Value = 69999.5
$ws.Range("L120").Value = 69999.5
$ws.Range("N120").Value = -77257.5

$ws.Range("H132").Value = 5528.1665
$ws.Range("I132").Value = 5454.3
$ws.Range("K132").Value = 16362.9
$ws.Range("M132").Value = -13832.9

$ws.Range("H134").Value = 3071.2104
$ws.Range("I134").Value = 2961.9412
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 8885.8236
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -6350.8236
$ws.Range("N134").Value = -17070

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 198082670
$ws.Range("I4").Value = 162516740
$ws.Range("K4").Value = 487550220
$ws.Range("M4").Value = -487550108

$ws.Range("H131").Value = 1809.2142
$ws.Range("I131").Value = 1777
$ws.Range("J131").Value = 1813.08
$ws.Range("K131").Value = 5331
$ws.Range("L131").Value = 5439.24
$ws.Range("M131").Value = -291
$ws.Range("N131").Value = -15519.24

$ws.Range("H132").Value = 559371.7
$ws.Range("J132").Value = 1115976.8
$ws.Range("L132").Value = 10043791.2
$ws.Range("N132").Value = -10048851.2

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H97").Value = 914.7368
$ws.Range("I97").Value = 743.3333
$ws.Range("K97").Value = 743.3333
$ws.Range("M97").Value = -247.3333

$ws.Range("H113").Value = 18524.666
$ws.Range("I113").Value = 2162.5
$ws.Range("J113").Value = 51249
$ws.Range("K113").Value = 2162.5
$ws.Range("L113").Value = 51249
$ws.Range("M113").Value = 7.5
$ws.Range("N113").Value = -55589

$ws.Range("H139").Value = 540000
$ws.Range("J139").Value = 540000
$ws.Range("L139").Value = 540000
$ws.Range("N139").Value = -550280

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H93").Value = 3280.3
$ws.Range("I93").Value = 1799.5
$ws.Range("J93").Value = 3650.5
$ws.Range("K93").Value = 1799.5
$ws.Range("L93").Value = 3650.5
$ws.Range("M93").Value = -551.5
$ws.Range("N93").Value = -6146.5

$ws.Range("H122").Value = 21882.611
$ws.Range("J122").Value = 18993
$ws.Range("L122").Value = 56979
$ws.Range("N122").Value = -61879

$ws.Range("H132").Value = 2505474.2
$ws.Range("I132").Value = 5003502
$ws.Range("J132").Value = 7446.5
$ws.Range("K132").Value = 15010506
$ws.Range("L132").Value = 22339.5
$ws.Range("M132").Value = -15007976
$ws.Range("N132").Value = -27399.5

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H14").Value = 83333750
$ws.Range("J14").Value = 1255
$ws.Range("L14").Value = 1255
$ws.Range("N14").Value = -1591

$ws.Range("H107").Value = 339.08334
$ws.Range("I107").Value = 277.35715
$ws.Range("K107").Value = 832.0714499999999
$ws.Range("M107").Value = 1087.92855

$ws.Range("H122").Value = 6282
$ws.Range("I122").Value = 6710.2
$ws.Range("K122").Value = 20130.6
$ws.Range("M122").Value = -17680.6

$ws.Range("H132").Value = 37915.43
$ws.Range("I132").Value = 39171.555
$ws.Range("K132").Value = 117514.665
$ws.Range("M132").Value = -114984.665
